$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:F3")
$range.NumberFormat = "@"

# Swap the runs/balls/fours/sixes values between row 2 and row 3 (columns C:F)
# Row 2 becomes what row 3 was (6, 4, 0, 1)
# Row 3 becomes what row 2 was (4, 1, 1, 0)
$ws.Range("C2").Value = "6"
$ws.Range("D2").Value = "4"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "1"

$ws.Range("C3").Value = "4"
$ws.Range("D3").Value = "1"
$ws.Range("E3").Value = "1"
$ws.Range("F3").Value = "0"
